$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 788,
# pushing all subsequent rows (788-826) down by one (to 789-827).
$ws.Rows.Item(788).Insert()

# Populate the newly inserted row 788 with the new record's data.
$ws.Cells.Item(788, 1).Value = 10
$ws.Cells.Item(788, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(788, 3).Value = 'La Araucanía'
$ws.Cells.Item(788, 4).Value = 45013
$ws.Cells.Item(788, 5).Value = 9
$ws.Cells.Item(788, 6).Value = 'Fruta'
$ws.Cells.Item(788, 7).Value = 100102
$ws.Cells.Item(788, 8).Value = 'Cítricos'
$ws.Cells.Item(788, 9).Value = 100102004
$ws.Cells.Item(788, 10).Value = 'Mandarina'
$ws.Cells.Item(788, 11).Value = 'Murcott'
$ws.Cells.Item(788, 12).Value = 'Primera'
$ws.Cells.Item(788, 13).Value = 55
$ws.Cells.Item(788, 14).Value = 18000
$ws.Cells.Item(788, 15).Value = 18000
$ws.Cells.Item(788, 16).Value = 18000
$ws.Cells.Item(788, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(788, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(788, 19).Value = 1000
$ws.Cells.Item(788, 20).Value = 18
